$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This applies the minimal set of cell-level changes between the previous
# revision and the new upload: two new scratch values in row 1, several
# course-code substitutions/shuffles in the Fall/Spring/Summer 2022 and
# 2023 blocks, and two new course rows (10 and 15) that were previously
# blank.

# --- New row 1 scratch values ---
$ws.Range("C1").Value = "dadaw"
# "32423" needs to land as text (a shared string), not the number 32423.
# Build it with TEXT() and paste-special the value back over the formula so
# the cell ends up a plain text cell (no quote-prefix/number-format noise).
$ws.Range("E1").Formula = "=TEXT(32423,""0"")"
$ws.Range("E1").Copy()
$ws.Range("E1").PasteSpecial(-4163)

# --- Fall 2022 column (A4:A9) course code reshuffle ---
$ws.Range("A4").Value = "PSYC 1101"
$ws.Range("C4").Value = "CPSC 3165"

$ws.Range("A5").Value = "POLS 1101"

$ws.Range("A6").Value = "PSYC 1105"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = "CPSC 4148"

$ws.Range("A7").Value = "DSCI 3111"
$ws.Range("C7").Value = "CPSC 4155"

$ws.Range("A8").Value = "CPSC 3121"
$ws.Range("C8").Value = "CPSC 4157"

# --- Fall 2023 column (A13:A15) reshuffle + new rows ---
$ws.Range("A13").Value = "CPSC 4175"
$ws.Range("A14").Value = "CPSC 4205"

# --- New row 10 (Fall 2022 block) ---
$ws.Range("A10").Value = "CYBR 4416"
$ws.Range("B10").Value = 1

# --- New row 15 (Fall 2023 block) ---
$ws.Range("A15").Value = "CPSC 4555"
$ws.Range("B15").Value = 3
